# Update to v3.0 x carriage
# - Row 2 (accuracy requirement): replace the Engineering Requirement with a
#   new statement and clear the Verification Process cell.
# - Append a new requirement row ("autonomous") to the Table1 table.
# - Move the active selection to D1 (and drop the old scrolled/selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "accurate" requirement gets reworded, Verification Process cleared ---
$ws.Cells.Item(2, 2).Value2 = "The system should be accurate to within +- 1% workable area"
$ws.Cells.Item(2, 3).ClearContents()
$ws.Rows.Item(2).RowHeight = 28.5

# --- Grow Table1 by one row, then fill it in like the row above it ---
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Range("A9:D9").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(10, 1).Value2 = "The system should be autonomous (no human intervention)"
$ws.Cells.Item(10, 2).Value2 = "At least for 8 hours (TODO - Look up official time periods for automation) "
$ws.Rows.Item(10).RowHeight = 42.75

# --- Update the view/selection ---
$ws.Range("D1").Select() | Out-Null
